# Expansão das análises automáticas: adiciona as colunas L, M e N
# (apoio_medio, contribuicoes, media_contribuicoes) à planilha existente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cabeçalhos (linha 1) -------------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Replica a formatação do cabeçalho existente (negrito, bordas, centralizado)
# para as novas células de cabeçalho.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Dados (linhas 2 a 7) --------------------------------------------------
$ws.Range("L2").Value = 92.07932629902824
$ws.Range("M2").Value = 228260
$ws.Range("N2").Value = 296.8270481144343

$ws.Range("L3").Value = 86.28489520037526
$ws.Range("M3").Value = 35293
$ws.Range("N3").Value = 578.5737704918033

$ws.Range("L4").Value = 88.79569566345432
$ws.Range("M4").Value = 180657
$ws.Range("N4").Value = 141.2486317435497

$ws.Range("L5").Value = 100.9425353902489
$ws.Range("M5").Value = 22989
$ws.Range("N5").Value = 221.0480769230769

$ws.Range("L6").Value = 19.4629156234702
$ws.Range("M6").Value = 2121
$ws.Range("N6").Value = 14.14

$ws.Range("L7").Value = 21.90935307361503
$ws.Range("M7").Value = 87
$ws.Range("N7").Value = 43.5
